$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.231.34"
$ws.Range("E2").Value = "  -0.24%  "

$ws.Range("D3").Value = "'1.859.07"
$ws.Range("E3").Value = "  -0.29%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.47%  "

$ws.Range("D5").Value = "'0.7103"
$ws.Range("E5").Value = "  +1.15%  "

$ws.Range("D6").Value = "'241.54"
$ws.Range("E6").Value = "  -0.45%  "

$ws.Range("E7").Value = "  -0.44%  "

$ws.Range("D8").Value = "'0.3095"
$ws.Range("E8").Value = "  -0.23%  "

$ws.Range("D9").Value = "'0.07739"
$ws.Range("E9").Value = "  -0.50%  "

$ws.Range("D10").Value = "'23.77"
$ws.Range("E10").Value = "  -1.69%  "

$ws.Range("D11").Value = "'0.07801"
$ws.Range("E11").Value = "  -2.72%  "

$ws.Range("D12").Value = "'1.863.23"
$ws.Range("E12").Value = "  +9.82%  "

$ws.Range("D13").Value = "'5.091"
$ws.Range("E13").Value = "  -1.38%  "

$ws.Range("D14").Value = "'91.98"
$ws.Range("E14").Value = "  -0.84%  "

$ws.Range("D15").Value = "'0.6863"
$ws.Range("E15").Value = "  -1.24%  "

$ws.Range("D16").Value = "'6.505"
$ws.Range("E16").Value = "  +2.60%  "

$ws.Range("D17").Value = "'0.000008414"
$ws.Range("E17").Value = "  +2.01%  "

$ws.Range("D18").Value = "'29.222.88"
$ws.Range("E18").Value = "  -0.33%  "

$ws.Range("D19").Value = "'250.58"
$ws.Range("E19").Value = "  +0.41%  "

$ws.Range("D20").Value = "'2.104.01"
$ws.Range("E20").Value = "  -0.86%  "

$ws.Range("D21").Value = "'12.82"
$ws.Range("E21").Value = "  -2.22%  "

$ws.Range("E22").Value = "  -0.48%  "

$ws.Range("D23").Value = "'7.506"
$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.45%  "

$ws.Range("D25").Value = "'0.1536"
$ws.Range("E25").Value = "  -1.07%  "

$ws.Range("D26").Value = "'160.05"
$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("D27").Value = "'8.841"
$ws.Range("E27").Value = "  -1.23%  "

$ws.Range("D28").Value = "'18.51"
$ws.Range("E28").Value = "  -0.28%  "

$ws.Range("D29").Value = "'1.562"
$ws.Range("E29").Value = "  +4.16%  "

$ws.Range("D30").Value = "'4.243"
$ws.Range("E30").Value = "  -0.40%  "

$ws.Range("D31").Value = "'4.224"
$ws.Range("E31").Value = "  -0.70%  "

$ws.Range("D32").Value = "'1.196"
$ws.Range("E32").Value = "  -1.63%  "

$ws.Range("D33").Value = "'0.05206"
$ws.Range("E33").Value = "  -0.59%  "

$ws.Range("D34").Value = "'0.7586"
$ws.Range("E34").Value = "  +2.44%  "

$ws.Range("D35").Value = "'1.841"
$ws.Range("E35").Value = "  -1.87%  "

$ws.Range("D36").Value = "'1.164"
$ws.Range("E36").Value = "  +0.58%  "

$ws.Range("D37").Value = "'2.709"
$ws.Range("E37").Value = "  -0.28%  "

$ws.Range("D38").Value = "'0.01858"
$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("D39").Value = "'1.215.02"
$ws.Range("E39").Value = "  -2.83%  "

$ws.Range("D40").Value = "'2.722"
$ws.Range("E40").Value = "  -0.67%  "

$ws.Range("D41").Value = "'0.8946"
$ws.Range("E41").Value = "  -0.34%  "

$ws.Range("D42").Value = "'109.80"
$ws.Range("E42").Value = "  -1.02%  "

$ws.Range("E43").Value = "  -0.40%  "

$ws.Range("D44").Value = "'5.593"
$ws.Range("E44").Value = "  -9.92%  "

$ws.Range("D45").Value = "'2.000.08"
$ws.Range("E45").Value = "  -1.43%  "

$ws.Range("D46").Value = "'0.5182"
$ws.Range("E46").Value = "  -0.49%  "

$ws.Range("D47").Value = "'64.44"
$ws.Range("E47").Value = "  -8.39%  "

$ws.Range("D48").Value = "'9.470"
$ws.Range("E48").Value = "  +1.77%  "

$ws.Range("D49").Value = "'0.00000000121"
$ws.Range("E49").Value = "  -3.88%  "

$ws.Range("D50").Value = "'1.748"
$ws.Range("E50").Value = "  -2.13%  "

$ws.Range("D51").Value = "'6.998"
$ws.Range("E51").Value = "  +0.46%  "
